$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "56.939.26"
Set-TextValue $ws.Range("E2") "  -1.68%  "
Set-TextValue $ws.Range("D3") "2.986.66"
Set-TextValue $ws.Range("D4") "0.999"
Set-TextValue $ws.Range("E4") "  -0.10%  "
Set-TextValue $ws.Range("D5") "499.31"
Set-TextValue $ws.Range("E5") "  -4.92%  "
Set-TextValue $ws.Range("D6") "138.11"
Set-TextValue $ws.Range("E6") "  -3.24%  "
Set-TextValue $ws.Range("E7") "  +0.05%  "
Set-TextValue $ws.Range("D8") "0.430"
Set-TextValue $ws.Range("E8") "  -3.44%  "
Set-TextValue $ws.Range("E9") "  -4.40%  "
Set-TextValue $ws.Range("D10") "0.108"
Set-TextValue $ws.Range("E10") "  -4.84%  "
Set-TextValue $ws.Range("D11") "0.358"
Set-TextValue $ws.Range("E11") "  -3.21%  "
Set-TextValue $ws.Range("D12") "3.495.71"
Set-TextValue $ws.Range("E12") "  -2.12%  "
Set-TextValue $ws.Range("E13") "  -2.36%  "
Set-TextValue $ws.Range("D14") "26.09"
Set-TextValue $ws.Range("E14") "  -3.08%  "
Set-TextValue $ws.Range("D15") "0.0000160"
Set-TextValue $ws.Range("E15") "  -7.30%  "
Set-TextValue $ws.Range("D16") "56.996.61"
Set-TextValue $ws.Range("E16") "  -1.51%  "
Set-TextValue $ws.Range("D17") "6.07"
Set-TextValue $ws.Range("E17") "  -2.26%  "
Set-TextValue $ws.Range("D18") "2.984.30"
Set-TextValue $ws.Range("E18") "  -2.21%  "
Set-TextValue $ws.Range("D19") "12.62"
Set-TextValue $ws.Range("E19") "  -2.44%  "
Set-TextValue $ws.Range("D20") "7.87"
Set-TextValue $ws.Range("E20") "  -4.02%  "
Set-TextValue $ws.Range("D21") "320.47"
Set-TextValue $ws.Range("E21") "  -6.24%  "
Set-TextValue $ws.Range("E22") "  -0.21%  "
Set-TextValue $ws.Range("D23") "5.70"
Set-TextValue $ws.Range("E23") "  -1.06%  "
Set-TextValue $ws.Range("D24") "0.491"
Set-TextValue $ws.Range("E24") "  -1.81%  "
Set-TextValue $ws.Range("D25") "63.63"
Set-TextValue $ws.Range("E25") "  -2.60%  "
Set-TextValue $ws.Range("E26") "  +0.26%  "
Set-TextValue $ws.Range("E27") "  -5.49%  "
Set-TextValue $ws.Range("D28") "0.0₃0893"
Set-TextValue $ws.Range("E28") "  -8.16%  "
Set-TextValue $ws.Range("D29") "6.53"
Set-TextValue $ws.Range("E29") "  -6.83%  "
Set-TextValue $ws.Range("D30") "7.06"
Set-TextValue $ws.Range("E30") "  -3.19%  "
Set-TextValue $ws.Range("E31") "  -5.31%  "
Set-TextValue $ws.Range("E32") "  -6.67%  "
Set-TextValue $ws.Range("D33") "20.22"
Set-TextValue $ws.Range("E33") "  -4.19%  "
Set-TextValue $ws.Range("D34") "155.99"
Set-TextValue $ws.Range("E34") "  -0.27%  "
Set-TextValue $ws.Range("D35") "4.56"
Set-TextValue $ws.Range("E35") "  -3.75%  "
Set-TextValue $ws.Range("D36") "5.79"
Set-TextValue $ws.Range("E36") "  -2.01%  "
Set-TextValue $ws.Range("E37") "  -6.83%  "
Set-TextValue $ws.Range("D38") "24.26"
Set-TextValue $ws.Range("E38") "  -6.19%  "
Set-TextValue $ws.Range("D39") "0.0664"
Set-TextValue $ws.Range("E39") "  -4.67%  "
Set-TextValue $ws.Range("B40") "OKB"
Set-TextValue $ws.Range("C40") "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue $ws.Range("D40") "37.72"
Set-TextValue $ws.Range("E40") "  +0.04%  "
Set-TextValue $ws.Range("B41") "RenzoRestakedETH"
Set-TextValue $ws.Range("C41") "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
Set-TextValue $ws.Range("D41") "3.011.99"
Set-TextValue $ws.Range("E41") "  -2.40%  "
Set-TextValue $ws.Range("D42") "0.999"
Set-TextValue $ws.Range("E42") "  -0.11%  "
Set-TextValue $ws.Range("D43") "3.74"
Set-TextValue $ws.Range("E43") "  -3.17%  "
Set-TextValue $ws.Range("D44") "0.642"
Set-TextValue $ws.Range("E44") "  -3.24%  "
Set-TextValue $ws.Range("D45") "2.204.90"
Set-TextValue $ws.Range("E45") "  -5.25%  "
Set-TextValue $ws.Range("E46") "  -6.17%  "
Set-TextValue $ws.Range("B47") "ONDO"
Set-TextValue $ws.Range("C47") "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
Set-TextValue $ws.Range("D47") "0.942"
Set-TextValue $ws.Range("E47") "  -8.84%  "
Set-TextValue $ws.Range("B48") "Cosmos"
Set-TextValue $ws.Range("C48") "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue $ws.Range("D48") "5.95"
Set-TextValue $ws.Range("E48") "  -1.59%  "
Set-TextValue $ws.Range("E49") "  -5.84%  "
Set-TextValue $ws.Range("D50") "19.28"
Set-TextValue $ws.Range("E50") "  -4.30%  "
Set-TextValue $ws.Range("D51") "1.80"
Set-TextValue $ws.Range("E51") "  -11.81%  "
